$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 12 and 15 (F:V content) ---
$ws.Cells.Item(12, 6).Value2 = 'GOAL FC'
$ws.Cells.Item(12, 7).Value2 = 3
$ws.Cells.Item(12, 8).Value2 = 'Red Star'
$ws.Cells.Item(12, 9).Value2 = 1
$ws.Cells.Item(12, 10).Value2 = 4.16
$ws.Cells.Item(12, 11).Value2 = '17/08/2023 11:39'
$ws.Cells.Item(12, 12).Value2 = 4.41
$ws.Cells.Item(12, 13).Value2 = '18/08/2023 19:27'
$ws.Cells.Item(12, 14).Value2 = 3.48
$ws.Cells.Item(12, 15).Value2 = '17/08/2023 11:39'
$ws.Cells.Item(12, 16).Value2 = 3.12
$ws.Cells.Item(12, 17).Value2 = '18/08/2023 19:27'
$ws.Cells.Item(12, 18).Value2 = 1.83
$ws.Cells.Item(12, 19).Value2 = '17/08/2023 11:39'
$ws.Cells.Item(12, 20).Value2 = 2.01
$ws.Cells.Item(12, 21).Value2 = '18/08/2023 19:27'
$ws.Cells.Item(12, 22).Value2 = 'https://www.betexplorer.com/football/france/national/goal-fc-red-star/hjeephhA/'
$ws.Cells.Item(15, 6).Value2 = 'Martigues'
$ws.Cells.Item(15, 7).Value2 = 2
$ws.Cells.Item(15, 8).Value2 = 'Versailles'
$ws.Cells.Item(15, 9).Value2 = 1
$ws.Cells.Item(15, 10).Value2 = 1.94
$ws.Cells.Item(15, 11).Value2 = '17/08/2023 11:39'
$ws.Cells.Item(15, 12).Value2 = 2.36
$ws.Cells.Item(15, 13).Value2 = '18/08/2023 19:22'
$ws.Cells.Item(15, 14).Value2 = 3.44
$ws.Cells.Item(15, 15).Value2 = '17/08/2023 11:39'
$ws.Cells.Item(15, 16).Value2 = 3.17
$ws.Cells.Item(15, 17).Value2 = '18/08/2023 19:22'
$ws.Cells.Item(15, 18).Value2 = 3.71
$ws.Cells.Item(15, 19).Value2 = '17/08/2023 11:39'
$ws.Cells.Item(15, 20).Value2 = 3.27
$ws.Cells.Item(15, 21).Value2 = '18/08/2023 19:29'
$ws.Cells.Item(15, 22).Value2 = 'https://www.betexplorer.com/football/france/national/martigues-versailles/IX2rmf8i/'

# --- Rotate rows 79-84 (F:V content): 79<-80, 80<-81, 81<-82, 82<-83, 83<-84, 84<-79 ---
$ws.Cells.Item(79, 6).Value2 = 'Chateauroux'
$ws.Cells.Item(79, 7).Value2 = 1
$ws.Cells.Item(79, 8).Value2 = 'Sochaux'
$ws.Cells.Item(79, 9).Value2 = 1
$ws.Cells.Item(79, 10).Value2 = 2.31
$ws.Cells.Item(79, 11).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(79, 12).Value2 = 2.88
$ws.Cells.Item(79, 13).Value2 = '11/10/2023 20:51'
$ws.Cells.Item(79, 14).Value2 = 2.92
$ws.Cells.Item(79, 15).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(79, 16).Value2 = 3.08
$ws.Cells.Item(79, 17).Value2 = '11/10/2023 20:51'
$ws.Cells.Item(79, 18).Value2 = 3.38
$ws.Cells.Item(79, 19).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(79, 20).Value2 = 2.69
$ws.Cells.Item(79, 21).Value2 = '11/10/2023 20:51'
$ws.Cells.Item(79, 22).Value2 = 'https://www.betexplorer.com/football/france/national/chateauroux-sochaux/IFH4C4ZA/'
$ws.Cells.Item(80, 6).Value2 = 'Cholet'
$ws.Cells.Item(80, 7).Value2 = 0
$ws.Cells.Item(80, 8).Value2 = 'Orleans'
$ws.Cells.Item(80, 9).Value2 = 1
$ws.Cells.Item(80, 10).Value2 = 2.44
$ws.Cells.Item(80, 11).Value2 = '10/10/2023 12:12'
$ws.Cells.Item(80, 12).Value2 = 3.01
$ws.Cells.Item(80, 13).Value2 = '11/10/2023 20:58'
$ws.Cells.Item(80, 14).Value2 = 3.14
$ws.Cells.Item(80, 15).Value2 = '10/10/2023 12:12'
$ws.Cells.Item(80, 16).Value2 = 3
$ws.Cells.Item(80, 17).Value2 = '11/10/2023 20:58'
$ws.Cells.Item(80, 18).Value2 = 3.07
$ws.Cells.Item(80, 19).Value2 = '10/10/2023 12:12'
$ws.Cells.Item(80, 20).Value2 = 2.64
$ws.Cells.Item(80, 21).Value2 = '11/10/2023 20:58'
$ws.Cells.Item(80, 22).Value2 = 'https://www.betexplorer.com/football/france/national/cholet-orleans/O61sw5LN/'
$ws.Cells.Item(81, 6).Value2 = 'Dijon'
$ws.Cells.Item(81, 7).Value2 = 4
$ws.Cells.Item(81, 8).Value2 = 'Epinal'
$ws.Cells.Item(81, 9).Value2 = 1
$ws.Cells.Item(81, 10).Value2 = 1.63
$ws.Cells.Item(81, 11).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(81, 12).Value2 = 1.38
$ws.Cells.Item(81, 13).Value2 = '11/10/2023 20:51'
$ws.Cells.Item(81, 14).Value2 = 3.56
$ws.Cells.Item(81, 15).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(81, 16).Value2 = 4.9
$ws.Cells.Item(81, 17).Value2 = '11/10/2023 20:59'
$ws.Cells.Item(81, 18).Value2 = 5.38
$ws.Cells.Item(81, 19).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(81, 20).Value2 = 8.35
$ws.Cells.Item(81, 21).Value2 = '11/10/2023 20:59'
$ws.Cells.Item(81, 22).Value2 = 'https://www.betexplorer.com/football/france/national/dijon-epinal/KE4QtmMb/'
$ws.Cells.Item(82, 6).Value2 = 'Marignane'
$ws.Cells.Item(82, 7).Value2 = 1
$ws.Cells.Item(82, 8).Value2 = 'Nimes'
$ws.Cells.Item(82, 9).Value2 = 1
$ws.Cells.Item(82, 10).Value2 = 2.78
$ws.Cells.Item(82, 11).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(82, 12).Value2 = 3.67
$ws.Cells.Item(82, 13).Value2 = '11/10/2023 20:57'
$ws.Cells.Item(82, 14).Value2 = 2.86
$ws.Cells.Item(82, 15).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(82, 16).Value2 = 2.96
$ws.Cells.Item(82, 17).Value2 = '11/10/2023 20:57'
$ws.Cells.Item(82, 18).Value2 = 2.78
$ws.Cells.Item(82, 19).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(82, 20).Value2 = 2.3
$ws.Cells.Item(82, 21).Value2 = '11/10/2023 20:57'
$ws.Cells.Item(82, 22).Value2 = 'https://www.betexplorer.com/football/france/national/marignane-nimes/A3cwvo6H/'
$ws.Cells.Item(83, 6).Value2 = 'Nancy'
$ws.Cells.Item(83, 7).Value2 = 1
$ws.Cells.Item(83, 8).Value2 = 'Red Star'
$ws.Cells.Item(83, 9).Value2 = 1
$ws.Cells.Item(83, 10).Value2 = 2.83
$ws.Cells.Item(83, 11).Value2 = '10/10/2023 12:12'
$ws.Cells.Item(83, 12).Value2 = 4.38
$ws.Cells.Item(83, 13).Value2 = '11/10/2023 20:56'
$ws.Cells.Item(83, 14).Value2 = 3.12
$ws.Cells.Item(83, 15).Value2 = '10/10/2023 12:12'
$ws.Cells.Item(83, 16).Value2 = 3.44
$ws.Cells.Item(83, 17).Value2 = '11/10/2023 20:56'
$ws.Cells.Item(83, 18).Value2 = 2.52
$ws.Cells.Item(83, 19).Value2 = '10/10/2023 12:12'
$ws.Cells.Item(83, 20).Value2 = 1.89
$ws.Cells.Item(83, 21).Value2 = '11/10/2023 20:56'
$ws.Cells.Item(83, 22).Value2 = 'https://www.betexplorer.com/football/france/national/nancy-red-star/n16Ir9in/'
$ws.Cells.Item(84, 6).Value2 = 'Niort'
$ws.Cells.Item(84, 7).Value2 = 1
$ws.Cells.Item(84, 8).Value2 = 'Versailles'
$ws.Cells.Item(84, 9).Value2 = 0
$ws.Cells.Item(84, 10).Value2 = 2.23
$ws.Cells.Item(84, 11).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(84, 12).Value2 = 1.89
$ws.Cells.Item(84, 13).Value2 = '11/10/2023 20:56'
$ws.Cells.Item(84, 14).Value2 = 3.03
$ws.Cells.Item(84, 15).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(84, 16).Value2 = 3.42
$ws.Cells.Item(84, 17).Value2 = '11/10/2023 20:56'
$ws.Cells.Item(84, 18).Value2 = 3.41
$ws.Cells.Item(84, 19).Value2 = '10/10/2023 12:20'
$ws.Cells.Item(84, 20).Value2 = 4.41
$ws.Cells.Item(84, 21).Value2 = '11/10/2023 20:56'
$ws.Cells.Item(84, 22).Value2 = 'https://www.betexplorer.com/football/france/national/niort-versailles/tO5MsT6h/'

# --- Add new row 94 (new match result appended to the bottom) ---
# Copy formatting (styles) from row 93 down to the new row 94 first,
# then populate the values explicitly.
$ws.Range("A93:V93").Copy()
$ws.Range("A94:V94").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(94, 1).Value2 = 93
$ws.Cells.Item(94, 2).Value2 = 'france'
$ws.Cells.Item(94, 3).Value2 = 'national'
$ws.Cells.Item(94, 4).Value2 = '2023-2024'
$ws.Cells.Item(94, 5).Value2 = 45224.79166666666
$ws.Cells.Item(94, 6).Value2 = 'Sochaux'
$ws.Cells.Item(94, 7).Value2 = 4
$ws.Cells.Item(94, 8).Value2 = 'Orleans'
$ws.Cells.Item(94, 9).Value2 = 2
$ws.Cells.Item(94, 10).Value2 = 2.22
$ws.Cells.Item(94, 11).Value2 = '25/10/2023 15:42'
$ws.Cells.Item(94, 12).Value2 = 2.28
$ws.Cells.Item(94, 13).Value2 = '25/10/2023 18:50'
$ws.Cells.Item(94, 14).Value2 = 3.06
$ws.Cells.Item(94, 15).Value2 = '25/10/2023 15:42'
$ws.Cells.Item(94, 16).Value2 = 3.11
$ws.Cells.Item(94, 17).Value2 = '25/10/2023 18:50'
$ws.Cells.Item(94, 18).Value2 = 3.39
$ws.Cells.Item(94, 19).Value2 = '25/10/2023 15:42'
$ws.Cells.Item(94, 20).Value2 = 3.52
$ws.Cells.Item(94, 21).Value2 = '25/10/2023 18:50'
$ws.Cells.Item(94, 22).Value2 = 'https://www.betexplorer.com/football/france/national/sochaux-orleans/4CiVc5KH/'
